$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D:D").Insert()

# Copy number formats/styles from the (shifted) former column D -- now column E -- into new column D
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Rows that never had a value in columns D:K before the insert should stay
# completely empty in the new column D (no stray formatted-but-empty cell)
$ws.Range("D5").Clear()
$ws.Range("D6").Clear()
$ws.Range("D37").Clear()
$ws.Range("D79").Clear()

# Populate new column D with the newly reported (most recent) year of data
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 1193100
$ws.Range("D9").Value = 747400
$ws.Range("D10").Value = 445700
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 1500
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 1160600
$ws.Range("D18").Value = 32500
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = 39900
$ws.Range("D22").Value = 13000
$ws.Range("D23").Value = 19500
$ws.Range("D24").Value = 1200
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 18300
$ws.Range("D27").Value = 18300
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = 18300
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 18300
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 11900
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 8250700
$ws.Range("D48").Value = 41000
$ws.Range("D49").Value = 47400
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 11031900
$ws.Range("D57").Value = 0
$ws.Range("D58").Value = "NA"
$ws.Range("D59").Value = 3045300
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 297700
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 9741300
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 1216600
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1290500
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 18300
$ws.Range("D83").Value = 7400
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 200900
$ws.Range("D91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -186500
$ws.Range("D96").Value = -46700
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -10100
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 4300
